$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 31.743396
$ws.Range("N2").Value = 95.230188
$ws.Range("O2").Value = 0.1189237443612096
$ws.Range("P2").Value = 0.1189237443612096
$ws.Range("Q2").Value = 0.15612460266
$ws.Range("R2").Value = 1.40512142394
$ws.Range("S2").Value = 0.1189237443612096
$ws.Range("T2").Value = 0.1189237443612096

$ws.Range("O3").Value = 0.1062760886263749
$ws.Range("P3").Value = 0.106276088626375
$ws.Range("S3").Value = 0.1062760886263749
$ws.Range("T3").Value = 0.106276088626375

$ws.Range("M4").Value = 70.51016133333333
$ws.Range("N4").Value = 211.530484
$ws.Range("O4").Value = 0.264159902780187
$ws.Range("P4").Value = 0.264159902780187
$ws.Range("Q4").Value = 0.3467924768244444
$ws.Range("R4").Value = 3.12113229142
$ws.Range("S4").Value = 0.264159902780187
$ws.Range("T4").Value = 0.264159902780187

$ws.Range("M5").Value = 6.595865666666666
$ws.Range("N5").Value = 19.787597
$ws.Range("O5").Value = 0.02471081047483217
$ws.Range("P5").Value = 0.02471081047483218
$ws.Range("Q5").Value = 0.03244066597055555
$ws.Range("R5").Value = 0.291965993735
$ws.Range("S5").Value = 0.02471081047483217
$ws.Range("T5").Value = 0.02471081047483218

$ws.Range("M6").Value = 64.15060166666667
$ws.Range("N6").Value = 192.451805
$ws.Range("O6").Value = 0.2403343912297365
$ws.Range("P6").Value = 0.2403343912297365
$ws.Range("Q6").Value = 0.3155140425305556
$ws.Range("R6").Value = 2.839626382775
$ws.Range("S6").Value = 0.2403343912297365
$ws.Range("T6").Value = 0.2403343912297365

$ws.Range("M7").Value = 65.55479199999999
$ws.Range("N7").Value = 196.664376
$ws.Range("O7").Value = 0.2455950625276598
$ws.Range("P7").Value = 0.2455950625276598
$ws.Range("Q7").Value = 0.3224203186533333
$ws.Range("R7").Value = 2.90178286788
$ws.Range("S7").Value = 0.2455950625276598
$ws.Range("T7").Value = 0.2455950625276598
